$wb = $excel.ActiveWorkbook

# --- Add the new "InvalidLogin" worksheet after the existing "ValidLogin" sheet ---
$validLogin = $wb.Worksheets.Item("ValidLogin")
$invalidLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# --- Populate InvalidLogin with header row + sample valid/invalid login rows ---
$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"
$invalidLogin.Range("A3").Value = "admin"
$invalidLogin.Range("B3").Value = "damager"

# --- Update the ValidLogin sheet's selection (was B3, now A2) and zoom to 220% ---
$validLogin.Select() | Out-Null
$validLogin.Range("A2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220

# --- InvalidLogin becomes the active/visible tab, selection B3, zoom 220% ---
$invalidLogin.Select() | Out-Null
$invalidLogin.Range("B3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220
